$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '61.039.15'
$ws.Range('D2').NumberFormat = "General"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.91%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.368.76'
$ws.Range('D3').NumberFormat = "General"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.92%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '405.72'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.33%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '133.92'
$ws.Range('D6').NumberFormat = "General"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +10.06%  '

$ws.Range('E7').Value = '  +2.54%  '

$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.674'
$ws.Range('D9').NumberFormat = "General"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +5.21%  '

$ws.Range('E10').Value = '  +7.47%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '42.39'
$ws.Range('D11').NumberFormat = "General"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.29%  '

$ws.Range('E12').Value = '  -0.80%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.885.61'
$ws.Range('D13').NumberFormat = "General"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.40%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '8.32'
$ws.Range('D14').NumberFormat = "General"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.71%  '

$ws.Range('E15').Value = '  +0.63%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.351.09'
$ws.Range('D16').NumberFormat = "General"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.56%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '60.995.21'
$ws.Range('D17').NumberFormat = "General"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.85%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.02'
$ws.Range('D18').NumberFormat = "General"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.41%  '

$ws.Range('E19').Value = '  +1.68%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0000127'
$ws.Range('D20').NumberFormat = "General"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +8.22%  '

$ws.Range('E21').Value = '  -3.15%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '84.29'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +10.53%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '312.61'
$ws.Range('D23').NumberFormat = "General"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.81%  '

$ws.Range('E24').Value = '  -0.84%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.12'
$ws.Range('D25').NumberFormat = "General"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.05%  '

$ws.Range('E26').Value = '  +11.78%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '29.38'
$ws.Range('D27').NumberFormat = "General"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.34%  '

$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.27'
$ws.Range('D28').NumberFormat = "General"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +8.19%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.52'
$ws.Range('D29').NumberFormat = "General"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -7.63%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.173'
$ws.Range('D30').NumberFormat = "General"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.16%  '

$ws.Range('E31').Value = '  +1.59%  '

$ws.Range('E32').Value = '  -0.03%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.27'
$ws.Range('D33').NumberFormat = "General"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.13%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '41.19'
$ws.Range('D34').NumberFormat = "General"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.82%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.47'
$ws.Range('D35').NumberFormat = "General"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.80%  '

$ws.Range('E36').Value = '  +0.16%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '51.85'
$ws.Range('D37').NumberFormat = "General"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.95%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.997'
$ws.Range('D38').NumberFormat = "General"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.27%  '

$ws.Range('E39').Value = '  -2.17%  '

$ws.Range('E40').Value = '  -4.53%  '

$ws.Range('E41').Value = '  +1.26%  '

$ws.Range('E42').Value = '  +2.30%  '

$ws.Range('E43').Value = '  +0.91%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.01'
$ws.Range('D44').NumberFormat = "General"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.22%  '

$ws.Range('E45').Value = '  +1.23%  '

$ws.Range('E46').Value = '  -4.00%  '

$ws.Range('E47').Value = '  +1.11%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '21.38'
$ws.Range('D48').NumberFormat = "General"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.10%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.118.93'
$ws.Range('D49').NumberFormat = "General"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.89%  '

$ws.Range('E50').Value = '  -3.58%  '

$ws.Range('E51').Value = '  -0.36%  '

